# Weekly data refresh: a new "Apio" price record for Macroferia Regional de
# Talca is inserted as the new row 150, pushing the existing rows 150-184
# down to 151-185 (mirrors the upstream weekly appender inserting the latest
# observation at the top of this sub-range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 150 (everything below shifts down by one row,
# inheriting row 149's formatting - matches the existing date-style column D).
$ws.Rows.Item(150).EntireRow.Insert()

$ws.Range("A150").Value = 5
$ws.Range("B150").Value = "Macroferia Regional de Talca"
$ws.Range("C150").Value = "Maule"
$ws.Range("D150").Value = "6/10/2022"
$ws.Range("E150").Value = 7
$ws.Range("F150").Value = 100112017
$ws.Range("G150").Value = "Apio"
$ws.Range("H150").Value = "Americana (o)"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 500
$ws.Range("K150").Value = 6000
$ws.Range("L150").Value = 6000
$ws.Range("M150").Value = 6000
$ws.Range("N150").Value = "$/docena de matas"
$ws.Range("O150").Value = "Provincia del Elquí"
$ws.Range("P150").Value = 1000
$ws.Range("Q150").Value = 6
$ws.Range("R150").Value = "Hortaliza"
